$wb = $excel.ActiveWorkbook

# Sheet "展览": F6 316 -> 317, F9 7926 -> 7931
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 317
$ws1.Range("F9").Value = 7931

# Sheet "全部类型": F6 316 -> 317, F11 7926 -> 7931
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 317
$ws4.Range("F11").Value = 7931
